# Updates with supplemental figures:
# refresh the Fixed-Residue PSSM values (columns B:K, rows 2:21) with the
# recomputed scores from the new supplemental-figure run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(-19.02726348975318, 2.098733218119495, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318),
  @(-19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, 1.281663195536741, -19.02726348975318, -19.02726348975318),
  @(-19.02726348975318, 1.986068044537355, -19.02726348975318, -19.02726348975318, 3.326608329179088, -19.02726348975318, 1.500306102456822, -19.02726348975318, 0.9550240791079935, -19.02726348975318),
  @(-19.02726348975318, 1.59883853827697, -19.02726348975318, -19.02726348975318, -19.02726348975318, 2.810649656062704, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318),
  @(-19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318),
  @(2.430805132485141, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318),
  @(-19.02726348975318, -19.02726348975318, -19.02726348975318, 2.001255725910853, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318),
  @(3.868689109149571, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318),
  @(-19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, 1.759992656064892, -19.02726348975318, 2.201133061460117),
  @(-19.02726348975318, -19.02726348975318, -19.02726348975318, 2.894940878540662, -19.02726348975318, 2.831003276652627, -19.02726348975318, -19.02726348975318, -19.02726348975318, 1.957486637068142),
  @(-19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318),
  @(-19.02726348975318, -19.02726348975318, -19.02726348975318, 2.353491899101568, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, 1.701368706530995, 1.755809079347581),
  @(-19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, 1.96034818969558),
  @(-19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318),
  @(-19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, 1.947930256179621, -19.02726348975318),
  @(-19.02726348975318, 1.813210962263946, 4.321925529685866, -19.02726348975318, -19.02726348975318, -19.02726348975318, 2.030313835346022, 2.074307663860516, 2.512395633834621, -19.02726348975318),
  @(-19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, 1.951997710502025, 2.007084801993433, 2.391050854364313, -19.02726348975318),
  @(-19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, -19.02726348975318, 1.664058325392062, 1.824156927609953, -19.02726348975318, -19.02726348975318),
  @(-19.02726348975318, 1.260630342370452, -19.02726348975318, -19.02726348975318, 3.317227751829304, -19.02726348975318, 1.688742893837619, 1.268856221452838, -19.02726348975318, 2.087125116675542),
  @(-19.02726348975318, 1.493570966648447, -19.02726348975318, 1.785657244259986, -19.02726348975318, 2.552936641553924, 1.497784109769844, -19.02726348975318, -19.02726348975318, -19.02726348975318),
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  for ($j = 0; $j -lt $row.Length; $j++) {
    $ws.Cells.Item($i + 2, $j + 2).Value = $row[$j]
  }
}

Write-Output "done"
